# Update cryptos list values (Price / Volume(1h)) for Mon Jul 10 09:54:35 UTC 2023 GitHub Actions run.
# Column D (Price) and column E (Volume 1h) are plain/formatted-text cells in this sheet
# (prices can look like "30.144.29" with a thousands separator, and percents carry
# leading/trailing padding spaces), so for any new Price value that LOOKS like a plain
# number (e.g. "0.9996"), we format the cell as Text ("@") before assigning it so Excel
# keeps it a literal string instead of auto-converting it to a numeric value, then reset
# the cell Style back to "Normal" so no stray number-format style sticks to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.144.29'
$ws.Range("D3").Value = '1.860.39'
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9996'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.92'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.45%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4694'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("E8").Value = '  -0.37%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2859'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06468'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.58'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.92%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07664'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.43%  '
$ws.Range("D13").Value = '1.856.50'
$ws.Range("E13").Value = '  -0.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '93.69'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.060'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.05%  '
$ws.Range("E16").Value = '  -0.88%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '268.16'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.27%  '
$ws.Range("D18").Value = '30.134.27'
$ws.Range("E18").Value = '  -0.61%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.28'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9999'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000007520'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.34%  '
$ws.Range("D22").Value = '2.102.05'
$ws.Range("E22").Value = '  -0.64%  '
$ws.Range("E23").Value = '  -0.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.152'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.093'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.01%  '
$ws.Range("E26").Value = '  -0.88%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '166.13'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.70'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.09%  '
$ws.Range("E29").Value = '  -3.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.376'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.58%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09812'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.72%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.499'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.207'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.988'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.86%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.04665'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.00%  '
$ws.Range("E36").Value = '  -2.51%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6837'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.40%  '
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01819'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.720'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.52%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.332'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.99%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '69.99'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9990'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8329'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.879'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.99%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.78'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.04%  '
$ws.Range("E47").Value = '  -3.00%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.216'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '925.40'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.57%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.908'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '34.15'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.86%  '
